# CISC480-Charter-GroupC.xlsx update
# - Tesing task: start date moved 10/20 -> 10/18, duration 3 -> 4 workdays
#   (adjusted-days formula E5 recalculates from 3 to 5)
# - Deployment task: start date moved 10/23 -> 10/24
#   (completion-date formula D6 recalculates accordingly)
# - Selection cursor moved from H4 to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Tesing"): Start Date (B5) and Duration (C5)
$ws.Range("B5").Value = 45948
$ws.Range("C5").Value = 4

# Row 6 ("Deployment"): Start Date (B6)
$ws.Range("B6").Value = 45954

# Move the active selection cell to B7
$null = $ws.Range("B7").Select()
